$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename "tfidf" prefix to "tf-idf" in header row 1 (columns W:AQ) ---
$ws.Range("W1").Value = "tf-idf_mean"
$ws.Range("X1").Value = "tf-idf_std"
$ws.Range("Y1").Value = "tf-idf_fold0"
$ws.Range("Z1").Value = "tf-idf_fold1"
$ws.Range("AA1").Value = "tf-idf_fold2"
$ws.Range("AB1").Value = "tf-idf_fold3"
$ws.Range("AC1").Value = "tf-idf_fold4"
$ws.Range("AD1").Value = "tf-idf_chi_mean"
$ws.Range("AE1").Value = "tf-idf_chi_std"
$ws.Range("AF1").Value = "tf-idf_chi_fold0"
$ws.Range("AG1").Value = "tf-idf_chi_fold1"
$ws.Range("AH1").Value = "tf-idf_chi_fold2"
$ws.Range("AI1").Value = "tf-idf_chi_fold3"
$ws.Range("AJ1").Value = "tf-idf_chi_fold4"
$ws.Range("AK1").Value = "tf-idf_pca_mean"
$ws.Range("AL1").Value = "tf-idf_pca_std"
$ws.Range("AM1").Value = "tf-idf_pca_fold0"
$ws.Range("AN1").Value = "tf-idf_pca_fold1"
$ws.Range("AO1").Value = "tf-idf_pca_fold2"
$ws.Range("AP1").Value = "tf-idf_pca_fold3"
$ws.Range("AQ1").Value = "tf-idf_pca_fold4"

# --- Update RF row (row 6) values ---
$ws.Range("B6").Value = 0.7992955421404832
$ws.Range("C6").Value = 0.02456917382021785
$ws.Range("D6").Value = 0.7759412684688143
$ws.Range("E6").Value = 0.7789183217979441
$ws.Range("F6").Value = 0.8179964190150948
$ws.Range("G6").Value = 0.7854046051389504
$ws.Range("H6").Value = 0.8382170962816125
$ws.Range("I6").Value = 0.8166684721988402
$ws.Range("J6").Value = 0.02393855611041602
$ws.Range("K6").Value = 0.7977082590310615
$ws.Range("L6").Value = 0.7884897360703812
$ws.Range("M6").Value = 0.8076054640195676
$ws.Range("N6").Value = 0.8397558015274496
$ws.Range("O6").Value = 0.8497831003457411
$ws.Range("P6").Value = 0.7266662225643733
$ws.Range("Q6").Value = 0.02464296665336834
$ws.Range("R6").Value = 0.7378955264593563
$ws.Range("S6").Value = 0.6995666599576723
$ws.Range("T6").Value = 0.7343827979103733
$ws.Range("U6").Value = 0.6984699434137116
$ws.Range("V6").Value = 0.7630161850807529
$ws.Range("W6").Value = 0.8123268668058288
$ws.Range("X6").Value = 0.02586530967901127
$ws.Range("Y6").Value = 0.7999305525282765
$ws.Range("Z6").Value = 0.798511299317751
$ws.Range("AA6").Value = 0.817884508831439
$ws.Range("AB6").Value = 0.7855302306915211
$ws.Range("AC6").Value = 0.8597777426601567
$ws.Range("AD6").Value = 0.8086470213208156
$ws.Range("AE6").Value = 0.03262951385801117
$ws.Range("AF6").Value = 0.7593811517844467
$ws.Range("AG6").Value = 0.8061713447363559
$ws.Range("AH6").Value = 0.8200065985401384
$ws.Range("AI6").Value = 0.7976468267577631
$ws.Range("AJ6").Value = 0.8600291847853738
$ws.Range("AK6").Value = 0.8165292199912176
$ws.Range("AL6").Value = 0.04884137377120463
$ws.Range("AM6").Value = 0.7616407528641572
$ws.Range("AN6").Value = 0.7633789204206198
$ws.Range("AO6").Value = 0.8598862347041327
$ws.Range("AP6").Value = 0.816327773182612
$ws.Range("AQ6").Value = 0.8814124187845668

# --- Update Ensemble row (row 7) values ---
$ws.Range("B7").Value = 0.8419413357459821
$ws.Range("C7").Value = 0.05147781709889671
$ws.Range("D7").Value = 0.8517184942716858
$ws.Range("F7").Value = 0.8812080691112947
$ws.Range("G7").Value = 0.8171924697660753
$ws.Range("H7").Value = 0.9029536771472256
$ws.Range("I7").Value = 0.8596024325052311
$ws.Range("J7").Value = 0.03192770070836615
$ws.Range("M7").Value = 0.8585770860399209
$ws.Range("P7").Value = 0.8316139843172679
$ws.Range("Q7").Value = 0.04567796631473559
$ws.Range("R7").Value = 0.8220589374004827
$ws.Range("T7").Value = 0.8401234329714804
$ws.Range("W7").Value = 0.8435653735638207
$ws.Range("X7").Value = 0.03945846908330625
$ws.Range("Z7").Value = 0.7755072252225952
$ws.Range("AA7").Value = 0.8824405782607986
$ws.Range("AB7").Value = 0.8491433324322363
$ws.Range("AC7").Value = 0.8811743529061017
$ws.Range("AD7").Value = 0.8750760578550937
$ws.Range("AE7").Value = 0.04104328605369267
$ws.Range("AG7").Value = 0.8055677619095177
$ws.Range("AH7").Value = 0.9034079839646323
$ws.Range("AI7").Value = 0.8587319843138891
$ws.Range("AJ7").Value = 0.9246596554109288
$ws.Range("AK7").Value = 0.8457190001481187
$ws.Range("AL7").Value = 0.03098750368584828
$ws.Range("AO7").Value = 0.8718124727633402
$ws.Range("AP7").Value = 0.8495894731634306
$ws.Range("AQ7").Value = 0.8820028739903869
